# AssureTestData.xlsx: repurpose the "contacts" sheet into a "LoginTest"
# login-test data sheet (title row + username/password/expected-result table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old sheet had 4 columns (title/firstname/lastname/company); the new
# layout only needs 3 (username/password/expected), so drop column D
# entirely (it disappears from the sheet, not just cleared).
$ws.Columns.Item(4).Delete()

# Make room for a new title row above the header row; this pushes the
# existing header row (now just 3 columns) and all data rows down by one.
$ws.Rows.Item(1).Insert()

# Row 2: header row (keeps the bold / yellow-fill style from the old row 1).
$ws.Range("A2").Value = "InputUsername"
$ws.Range("B2").Value = "InputPassword"
$ws.Range("C2").Value = "ExpectedUsername"

# Row 3: first data row.
$ws.Range("A3").Value = "abhilasha.jha@northgateps.com"
$ws.Range("B3").Value = "N0rthg4t31"
$ws.Range("C3").Value = "Hello, Abhilashaa"

# Row 4: second data row (same values as row 3).
$ws.Range("A4").Value = "abhilasha.jha@northgateps.com"
$ws.Range("B4").Value = "N0rthg4t31"
$ws.Range("C4").Value = "Hello, Abhilashaa"

# Row 1: sheet title, filled in after the table body.
$ws.Range("A1").Value = "LoginTest"

# Row 5: third data row (uses a variant/incorrect password).
$ws.Range("A5").Value = "abhilasha.jha@northgateps.com"
$ws.Range("B5").Value = "N0rthg4t311"
$ws.Range("C5").Value = "Hello, Abhilashaa"

# Resize the columns to fit the new content (best-fit by character width).
$ws.Columns.Item(1).ColumnWidth = 26.833333333333332
$ws.Columns.Item(2).ColumnWidth = 13
$ws.Columns.Item(3).ColumnWidth = 16.5

# Put the selection back on the top-left cell.
$ws.Range("A1").Select()

# Rename the sheet to match its new purpose.
$ws.Name = "LoginTest"
